$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# F4 currently holds "low" (comparing primary key tolerance/value); change it to "High"
$ws.Range("F4").Value = "High"
